$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 279.33334
$ws.Range("I2").Value = 78.833336
$ws.Range("K2").Value = 78.833336
$ws.Range("M2").Value = 34.166664

$ws.Range("H15").Value = 888.7442
$ws.Range("I15").Value = 888.7442
$ws.Range("K15").Value = 2666.2326
$ws.Range("M15").Value = -2497.2326

$ws.Range("H17").Value = 1883.1578
$ws.Range("J17").Value = 1883.1578
$ws.Range("L17").Value = 5649.4734
$ws.Range("N17").Value = -5985.4734

$ws.Range("H19").Value = 1170.7142
$ws.Range("I19").Value = 1479
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 1479
$ws.Range("L19").Value = 400
$ws.Range("M19").Value = -1304
$ws.Range("N19").Value = -750

$ws.Range("H21").Value = 10850
$ws.Range("I21").Value = 1700
$ws.Range("K21").Value = 1700
$ws.Range("M21").Value = -1232

$ws.Range("H23").Value = 10850
$ws.Range("I23").Value = 1700
$ws.Range("K23").Value = 1700
$ws.Range("M23").Value = -1466

$ws.Range("H33").Value = 251.28572
$ws.Range("I33").Value = 218.16667
$ws.Range("J33").Value = 450
$ws.Range("K33").Value = 218.16667
$ws.Range("L33").Value = 450
$ws.Range("M33").Value = 10.83332999999999
$ws.Range("N33").Value = -908

$ws.Range("H41").Value = 72065.57000000001
$ws.Range("I41").Value = 83.09090999999999
$ws.Range("K41").Value = 83.09090999999999
$ws.Range("M41").Value = 356.90909

$ws.Range("H62").Value = 10962.8125
$ws.Range("I62").Value = 10633.25
$ws.Range("J62").Value = 11951.5
$ws.Range("K62").Value = 10633.25
$ws.Range("L62").Value = 11951.5
$ws.Range("M62").Value = -10009.25
$ws.Range("N62").Value = -13199.5

$ws.Range("H65").Value = 10962.8125
$ws.Range("I65").Value = 10633.25
$ws.Range("J65").Value = 11951.5
$ws.Range("K65").Value = 53166.25
$ws.Range("L65").Value = 59757.5
$ws.Range("M65").Value = -50046.25
$ws.Range("N65").Value = -65997.5

$ws.Range("H87").Value = 65123.19
$ws.Range("J87").Value = 70199.28
$ws.Range("L87").Value = 70199.28
$ws.Range("N87").Value = -72695.28

$ws.Range("H90").Value = 65123.19
$ws.Range("J90").Value = 70199.28
$ws.Range("L90").Value = 210597.84
$ws.Range("N90").Value = -223077.84

$ws.Range("H97").Value = 2242.7778
$ws.Range("J97").Value = 4072
$ws.Range("L97").Value = 12216
$ws.Range("N97").Value = -13208

$ws.Range("H112").Value = 1329.1212
$ws.Range("J112").Value = 1358.0968
$ws.Range("L112").Value = 4074.2904
$ws.Range("N112").Value = -6290.2904

$ws.Range("H132").Value = 2833.6
$ws.Range("I132").Value = 2772.2632
$ws.Range("K132").Value = 8316.7896
$ws.Range("M132").Value = -5786.7896

$ws.Range("H133").Value = 80277.5
$ws.Range("J133").Value = 80277.5
$ws.Range("L133").Value = 80277.5
$ws.Range("N133").Value = -90397.5

$ws.Range("H136").Value = 99331
$ws.Range("J136").Value = 99331
$ws.Range("L136").Value = 99331
$ws.Range("N136").Value = -109531

$ws.Range("H137").Value = 2309.8235
$ws.Range("I137").Value = 1306
$ws.Range("J137").Value = 2728.0833
$ws.Range("K137").Value = 3918
$ws.Range("L137").Value = 8184.249899999999
$ws.Range("M137").Value = -1368
$ws.Range("N137").Value = -13284.2499

$ws.Range("H138").Value = 1853.9807
$ws.Range("J138").Value = 3167.5264
$ws.Range("L138").Value = 9502.5792
$ws.Range("N138").Value = -19782.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3254.3208
$ws.Range("I32").Value = 3385.0833
$ws.Range("J32").Value = 1999
$ws.Range("K32").Value = 3385.0833
$ws.Range("L32").Value = 1999
$ws.Range("M32").Value = -3098.0833
$ws.Range("N32").Value = -2573

$ws.Range("H45").Value = 2690.2727
$ws.Range("I45").Value = 3177.0715
$ws.Range("K45").Value = 3177.0715
$ws.Range("M45").Value = -2800.0715

$ws.Range("H64").Value = 71199.39999999999
$ws.Range("J64").Value = 72666
$ws.Range("L64").Value = 72666
$ws.Range("N64").Value = -73162

$ws.Range("H67").Value = 71199.39999999999
$ws.Range("J67").Value = 72666
$ws.Range("L67").Value = 72666
$ws.Range("N67").Value = -74382

$ws.Range("H102").Value = 4731.5
$ws.Range("I102").Value = 3920.889
$ws.Range("K102").Value = 3920.889
$ws.Range("M102").Value = -2298.889

$ws.Range("H132").Value = 3625.95
$ws.Range("I132").Value = 3667.1538
$ws.Range("K132").Value = 11001.4614
$ws.Range("M132").Value = -8471.4614

$ws.Range("H134").Value = 89000
$ws.Range("J134").Value = 89000
$ws.Range("L134").Value = 89000
$ws.Range("N134").Value = -99140

$ws.Range("H139").Value = 88569.71000000001
$ws.Range("J139").Value = 88569.71000000001
$ws.Range("L139").Value = 88569.71000000001
$ws.Range("N139").Value = -98849.71000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 472
$ws.Range("I7").Value = 101
$ws.Range("K7").Value = 101
$ws.Range("M7").Value = 12

$ws.Range("H86").Value = 1467.3
$ws.Range("I86").Value = 1459.125
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1459.125
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -336.125
$ws.Range("N86").Value = -3746

$ws.Range("H89").Value = 1467.3
$ws.Range("I89").Value = 1459.125
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 7295.625
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -1679.625
$ws.Range("N89").Value = -18732

$ws.Range("H122").Value = 123333.336
$ws.Range("J122").Value = 123333.336
$ws.Range("L122").Value = 123333.336
$ws.Range("N122").Value = -133133.336

$ws.Range("H124").Value = 44000
$ws.Range("J124").Value = 44000
$ws.Range("L124").Value = 44000
$ws.Range("N124").Value = -53820

$ws.Range("H125").Value = 149999
$ws.Range("J125").Value = 149999
$ws.Range("L125").Value = 149999
$ws.Range("N125").Value = -159839

$ws.Range("H140").Value = 136662.33
$ws.Range("J140").Value = 136662.33
$ws.Range("L140").Value = 136662.33
$ws.Range("N140").Value = -147022.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 3329.8
$ws.Range("J12").Value = 4750
$ws.Range("L12").Value = 4750
$ws.Range("N12").Value = -5090

$ws.Range("H68").Value = 69519.55
$ws.Range("J68").Value = 69519.55
$ws.Range("L68").Value = 69519.55
$ws.Range("N68").Value = -71017.55

$ws.Range("H71").Value = 69519.55
$ws.Range("J71").Value = 69519.55
$ws.Range("L71").Value = 208558.65
$ws.Range("N71").Value = -216046.65

$ws.Range("H107").Value = 3773.5
$ws.Range("I107").Value = 3773.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3773.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1853.5
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 130458664
$ws.Range("I4").Value = 150118430
$ws.Range("K4").Value = 450355290
$ws.Range("M4").Value = -450355178

$ws.Range("H14").Value = 143170
$ws.Range("I14").Value = 143170
$ws.Range("K14").Value = 429510
$ws.Range("M14").Value = -429337

$ws.Range("H48").Value = 874.5
$ws.Range("J48").Value = 1499
$ws.Range("L48").Value = 4497
$ws.Range("N48").Value = -4997

$ws.Range("H98").Value = 363.6154
$ws.Range("J98").Value = 348.44446
$ws.Range("L98").Value = 1045.33338
$ws.Range("N98").Value = -4041.33338

$ws.Range("H105").Value = 17845.666
$ws.Range("J105").Value = 17845.666
$ws.Range("L105").Value = 53536.99800000001
$ws.Range("N105").Value = -58778.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 142262.4
$ws.Range("J82").Value = 140328
$ws.Range("L82").Value = 140328
$ws.Range("N82").Value = -141094

$ws.Range("H85").Value = 142262.4
$ws.Range("J85").Value = 140328
$ws.Range("L85").Value = 140328
$ws.Range("N85").Value = -142980

$ws.Range("H122").Value = 3853.0425
$ws.Range("I122").Value = 2108
$ws.Range("K122").Value = 6324
$ws.Range("M122").Value = -3874

$ws.Range("H126").Value = 2600.84
$ws.Range("I126").Value = 2333.875
$ws.Range("J126").Value = 2726.4707
$ws.Range("K126").Value = 7001.625
$ws.Range("L126").Value = 8179.4121
$ws.Range("M126").Value = -4531.625
$ws.Range("N126").Value = -13119.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 34801
$ws.Range("J101").Value = 34801
$ws.Range("L101").Value = 34801
$ws.Range("N101").Value = -41291

$ws.Range("H113").Value = 1682.6285
$ws.Range("I113").Value = 843.7917
$ws.Range("J113").Value = 3512.818
$ws.Range("K113").Value = 2531.3751
$ws.Range("L113").Value = 10538.454
$ws.Range("M113").Value = -361.3751000000002
$ws.Range("N113").Value = -14878.454

$ws.Range("H132").Value = 2665.796
$ws.Range("I132").Value = 1837.1892
$ws.Range("K132").Value = 5511.5676
$ws.Range("M132").Value = -2981.5676

$ws.Range("H136").Value = 18181536
$ws.Range("I136").Value = 21589324
$ws.Range("K136").Value = 64767972
$ws.Range("M136").Value = -64765422
